# Rebuild paragraph 1 runs with spell-check markers around "datetime" and "picker"
$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range.Duplicate
$p1.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">W </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>datetime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>picker</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wyświetlaj data dzisiejsza</w:t></w:r><w:r><w:t xml:space="preserve"> podczas dodawania nowego wniosku</w:t></w:r></w:p>')

# Rebuild paragraph 2 runs with spell-check markers around the English words
$p2 = $d.Paragraphs(2).Range.Duplicate
$p2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Object </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>reference</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> not set to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>an</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>instance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>an</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>object</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> co to znaczy</w:t></w:r></w:p>')

# Rebuild paragraph 3 runs with spell-check markers around "Negat"
$p3 = $d.Paragraphs(3).Range.Duplicate
$p3.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Dodać status wniosku: Złożony, Oferta, Decyzja, Wypłacony, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Negat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Rezygnacja Klienta</w:t></w:r></w:p>')

# Append a new empty paragraph at the very end of the document body (after
# paragraph 4, before the sectPr) - collapse the whole-document range to its
# end so the insertion lands after the last paragraph instead of replacing it.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
